$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose values look numeric,
# so Excel stores them as text (matching original inlineStr text cells).
$textPriceRows = @(5, 6, 7, 9, 10, 11, 12, 13, 14, 17, 20, 21, 23, 26, 28, 29, 31, 32, 33, 35, 38, 42, 43, 45, 47, 49, 50)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item([int]$r, 4).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '52.303.44'
$ws.Cells.Item(2, 5).Value = '  +1.91%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.798.68'
$ws.Cells.Item(3, 5).Value = '  +1.53%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '346.13'
$ws.Cells.Item(5, 5).Value = '  +4.70%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '116.80'
$ws.Cells.Item(6, 5).Value = '  +1.12%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.553'
$ws.Cells.Item(7, 5).Value = '  +4.07%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.05%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.586'
$ws.Cells.Item(9, 5).Value = '  +2.34%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '43.32'
$ws.Cells.Item(10, 5).Value = '  +4.50%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.0857'
$ws.Cells.Item(11, 5).Value = '  +3.40%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '20.12'
$ws.Cells.Item(12, 5).Value = '  -0.82%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '0.131'
$ws.Cells.Item(13, 5).Value = '  +1.76%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '7.83'
$ws.Cells.Item(14, 5).Value = '  +2.27%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.235.95'
$ws.Cells.Item(15, 5).Value = '  +1.51%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.799.74'
$ws.Cells.Item(16, 5).Value = '  +2.00%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '0.891'
$ws.Cells.Item(17, 5).Value = '  +0.34%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '52.202.17'
$ws.Cells.Item(18, 5).Value = '  +1.74%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +6.41%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '7.12'
$ws.Cells.Item(20, 5).Value = '  +3.87%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '13.43'
$ws.Cells.Item(21, 5).Value = '  -1.33%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +1.90%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '70.27'
$ws.Cells.Item(23, 5).Value = '  +0.14%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -3.61%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +6.34%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '26.71'
$ws.Cells.Item(26, 5).Value = '  -0.75%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.00%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '10.28'
$ws.Cells.Item(28, 5).Value = '  -0.72%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '2.24'
$ws.Cells.Item(29, 5).Value = '  +0.63%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.53%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '35.07'
$ws.Cells.Item(31, 5).Value = '  -1.82%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '50.29'
$ws.Cells.Item(32, 5).Value = '  +0.25%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '5.72'
$ws.Cells.Item(33, 5).Value = '  +1.56%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.08%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.0412'
$ws.Cells.Item(35, 5).Value = '  +16.84%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.68%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.03%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '18.91'
$ws.Cells.Item(38, 5).Value = '  -2.53%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.75%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -0.28%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +21.43%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).Value = '128.39'
$ws.Cells.Item(42, 5).Value = '  -0.81%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = '23.55'
$ws.Cells.Item(43, 5).Value = '  -0.81%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +2.31%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.32'
$ws.Cells.Item(45, 5).Value = '  +0.56%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -2.06%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(47, 4).Value = '2.39'
$ws.Cells.Item(47, 5).Value = '  +5.98%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).Value = '2.074.93'
$ws.Cells.Item(48, 5).Value = '  -2.01%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.987'
$ws.Cells.Item(49, 5).Value = '  +17.93%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '5.54'
$ws.Cells.Item(50, 5).Value = '  -0.10%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -1.35%  '
